$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2354.75
$ws.Range("I138").Value = 2135.1
$ws.Range("J138").Value = 2476.7778
$ws.Range("K138").Value = 6405.299999999999
$ws.Range("L138").Value = 7430.3334
$ws.Range("M138").Value = -1265.299999999999
$ws.Range("N138").Value = -17710.3334
$ws.Range("H141").Value = 6576
$ws.Range("I141").Value = 2772.5
$ws.Range("J141").Value = 11647.333
$ws.Range("K141").Value = 8317.5
$ws.Range("L141").Value = 34941.999
$ws.Range("M141").Value = -3137.5
$ws.Range("N141").Value = -45301.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 362268.22
$ws.Range("I32").Value = 416048
$ws.Range("K32").Value = 416048
$ws.Range("M32").Value = -415761
$ws.Range("H45").Value = 3600.875
$ws.Range("I45").Value = 2787.3333
$ws.Range("J45").Value = 4089
$ws.Range("K45").Value = 2787.3333
$ws.Range("L45").Value = 4089
$ws.Range("M45").Value = -2410.3333
$ws.Range("N45").Value = -4843
$ws.Range("H74").Value = 1559.5
$ws.Range("J74").Value = 1470.3
$ws.Range("L74").Value = 1470.3
$ws.Range("N74").Value = -3218.3
$ws.Range("H76").Value = 30288
$ws.Range("J76").Value = 30288
$ws.Range("L76").Value = 30288
$ws.Range("N76").Value = -30964
$ws.Range("H77").Value = 1559.5
$ws.Range("J77").Value = 1470.3
$ws.Range("L77").Value = 7351.5
$ws.Range("N77").Value = -16087.5
$ws.Range("H79").Value = 30288
$ws.Range("J79").Value = 30288
$ws.Range("L79").Value = 30288
$ws.Range("N79").Value = -32628
$ws.Range("H97").Value = 996.2619
$ws.Range("I97").Value = 972
$ws.Range("J97").Value = 1064.6364
$ws.Range("K97").Value = 972
$ws.Range("L97").Value = 1064.6364
$ws.Range("M97").Value = -476
$ws.Range("N97").Value = -2056.6364
$ws.Range("H132").Value = 9987.333000000001
$ws.Range("I132").Value = 18302
$ws.Range("J132").Value = 4444.222
$ws.Range("K132").Value = 54906
$ws.Range("L132").Value = 13332.666
$ws.Range("M132").Value = -52376
$ws.Range("N132").Value = -18392.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1926.6774
$ws.Range("I20").Value = 1594.4445
$ws.Range("J20").Value = 2386.6924
$ws.Range("K20").Value = 1594.4445
$ws.Range("L20").Value = 2386.6924
$ws.Range("M20").Value = -1347.4445
$ws.Range("N20").Value = -2880.6924
$ws.Range("H107").Value = 1466.6666
$ws.Range("I107").Value = 933.3333
$ws.Range("K107").Value = 933.3333
$ws.Range("M107").Value = 986.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3940.151
$ws.Range("I31").Value = 855.4583
$ws.Range("J31").Value = 6493
$ws.Range("K31").Value = 855.4583
$ws.Range("L31").Value = 6493
$ws.Range("M31").Value = -560.4583
$ws.Range("N31").Value = -7083
$ws.Range("H34").Value = 3940.151
$ws.Range("I34").Value = 855.4583
$ws.Range("J34").Value = 6493
$ws.Range("K34").Value = 855.4583
$ws.Range("L34").Value = 6493
$ws.Range("M34").Value = -653.4583
$ws.Range("N34").Value = -6897
$ws.Range("H58").Value = 2188.56
$ws.Range("I58").Value = 2340
$ws.Range("J58").Value = 1582.8
$ws.Range("K58").Value = 2340
$ws.Range("L58").Value = 1582.8
$ws.Range("M58").Value = -2137
$ws.Range("N58").Value = -1988.8
$ws.Range("H63").Value = 100000
$ws.Range("J63").Value = 100000
$ws.Range("L63").Value = 100000
$ws.Range("N63").Value = -101372
$ws.Range("H66").Value = 100000
$ws.Range("J66").Value = 100000
$ws.Range("L66").Value = 300000
$ws.Range("N66").Value = -306864
$ws.Range("H136").Value = 2188.56
$ws.Range("I136").Value = 2340
$ws.Range("J136").Value = 1582.8
$ws.Range("K136").Value = 7020
$ws.Range("L136").Value = 4748.4
$ws.Range("M136").Value = -4470
$ws.Range("N136").Value = -9848.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 563.5091
$ws.Range("I5").Value = 436.9737
$ws.Range("J5").Value = 846.35297
$ws.Range("K5").Value = 1310.9211
$ws.Range("L5").Value = 2539.05891
$ws.Range("M5").Value = -1198.9211
$ws.Range("N5").Value = -2763.05891
$ws.Range("H68").Value = 1287.1632
$ws.Range("I68").Value = 698.4
$ws.Range("J68").Value = 1488.7946
$ws.Range("K68").Value = 2095.2
$ws.Range("L68").Value = 4466.3838
$ws.Range("M68").Value = -1284.2
$ws.Range("N68").Value = -6088.3838
$ws.Range("H71").Value = 1287.1632
$ws.Range("I71").Value = 698.4
$ws.Range("J71").Value = 1488.7946
$ws.Range("K71").Value = 6285.599999999999
$ws.Range("L71").Value = 13399.1514
$ws.Range("M71").Value = -2229.599999999999
$ws.Range("N71").Value = -21511.1514
$ws.Range("H113").Value = 782.9545000000001
$ws.Range("I113").Value = 506.80768
$ws.Range("J113").Value = 1181.8334
$ws.Range("K113").Value = 1520.42304
$ws.Range("L113").Value = 3545.5002
$ws.Range("M113").Value = 649.5769599999999
$ws.Range("N113").Value = -7885.5002
$ws.Range("H131").Value = 1174.2609
$ws.Range("J131").Value = 1248.75
$ws.Range("L131").Value = 3746.25
$ws.Range("N131").Value = -13826.25
$ws.Range("H135").Value = 563.5091
$ws.Range("I135").Value = 436.9737
$ws.Range("J135").Value = 846.35297
$ws.Range("K135").Value = 3932.7633
$ws.Range("L135").Value = 7617.17673
$ws.Range("M135").Value = -1397.7633
$ws.Range("N135").Value = -12687.17673
$ws.Range("H138").Value = 1857.6666
$ws.Range("I138").Value = 727
$ws.Range("K138").Value = 2181
$ws.Range("M138").Value = 2959
$ws.Range("H140").Value = 1790.4839
$ws.Range("I140").Value = 1223.8
$ws.Range("K140").Value = 3671.4
$ws.Range("M140").Value = 1508.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3757.3542
$ws.Range("I132").Value = 3459.7097
$ws.Range("J132").Value = 4300.1177
$ws.Range("K132").Value = 10379.1291
$ws.Range("L132").Value = 12900.3531
$ws.Range("M132").Value = -7849.1291
$ws.Range("N132").Value = -17960.3531
$ws.Range("H141").Value = 43765.453
$ws.Range("J141").Value = 43765.453
$ws.Range("L141").Value = 43765.453
$ws.Range("N141").Value = -54125.453

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 24676.5
$ws.Range("J123").Value = 24676.5
$ws.Range("L123").Value = 24676.5
$ws.Range("N123").Value = -34476.5
$ws.Range("H136").Value = 3026.375
$ws.Range("I136").Value = 3245.75
$ws.Range("J136").Value = 2807
$ws.Range("K136").Value = 9737.25
$ws.Range("L136").Value = 8421
$ws.Range("M136").Value = -7187.25
$ws.Range("N136").Value = -13521
